# Add Nordic (Leirubakki) rating column to Evans2010 pole list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell L1, styled like the other header cells (copy format from K1)
$ws.Range("L1").Value = "Leirubakki_rating"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: L2:L3 = "B", L4:L15 = "A"
$ws.Range("L2").Value = "B"
$ws.Range("L3").Value = "B"
$ws.Range("L4:L15").Value = "A"

# Row 2 (header wrap row) grows a bit taller to fit the new column text
$ws.Rows(2).RowHeight = 85

# Move the active selection to the new bottom-right corner, matching the
# author's last selection position after adding the column
$ws.Range("L16").Select()
